# Actualización automática del tracker
# Adds two new result rows (122 y 123) al final de la hoja del tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = 14807182; B = "2025-10-08"; C = "Mats Rosenkranz"; D = "Micah Braswell";   E = "Gana Micah Braswell"; F = 1.73 },
    @{ A = 14807076; B = "2025-10-09"; C = "Radu Albot";      D = "Dimitar Kuzmanov"; E = "Gana Radu Albot";     F = 2.5 }
)

$newRow = 122

foreach ($r in $rows) {
    # Column A: numeric event id
    $ws.Cells.Item($newRow, 1).Value = $r.A

    # Column B: keep the date as literal text (matches existing rows, which
    # store "fecha" as inline text rather than a real Excel date). A leading
    # apostrophe forces Excel to treat it as a text literal instead of
    # converting it into a date serial number.
    $ws.Cells.Item($newRow, 2).Formula = "'" + $r.B
    $ws.Cells.Item($newRow, 2).Style = "Normal"

    # Columns C, D, E: plain text
    $ws.Cells.Item($newRow, 3).Value = $r.C
    $ws.Cells.Item($newRow, 4).Value = $r.D
    $ws.Cells.Item($newRow, 5).Value = $r.E

    # Column F: numeric odds
    $ws.Cells.Item($newRow, 6).Value = $r.F

    # Columns G, H: still pending result, stored as an explicit empty
    # string (same shape as the rows above, not simply a blank cell).
    $ws.Cells.Item($newRow, 7).Formula = "'"
    $ws.Cells.Item($newRow, 7).Style = "Normal"
    $ws.Cells.Item($newRow, 8).Formula = "'"
    $ws.Cells.Item($newRow, 8).Style = "Normal"

    $newRow = $newRow + 1
}
